$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.417.43'
$ws.Range('E2').Value = '  -0.39%  '
$ws.Range('D3').Value = '3.320.23'
$ws.Range('E3').Value = '  -0.06%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '586.66'
$ws.Range('E5').Value = '  +2.36%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '182.82'
$ws.Range('E6').Value = '  +0.28%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.646'
$ws.Range('E7').Value = '  +7.99%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  -2.01%  '
$ws.Range('E10').Value = '  +2.35%  '
$ws.Range('E11').Value = '  -0.20%  '
$ws.Range('D12').Value = '3.897.99'
$ws.Range('E12').Value = '  -0.08%  '
$ws.Range('E13').Value = '  -4.34%  '
$ws.Range('D14').Value = '66.425.63'
$ws.Range('E14').Value = '  -0.49%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.42'
$ws.Range('E15').Value = '  -2.76%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '3.322.77'
$ws.Range('E16').Value = '  +0.78%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000164'
$ws.Range('E17').Value = '  -2.09%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '431.55'
$ws.Range('E18').Value = '  -0.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.33'
$ws.Range('E19').Value = '  -2.31%  '
$ws.Range('E20').Value = '  -2.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.45'
$ws.Range('E21').Value = '  -2.49%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '72.30'
$ws.Range('E22').Value = '  -1.68%  '
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.71'
$ws.Range('E24').Value = '  +0.73%  '
$ws.Range('D25').Value = '3.440.63'
$ws.Range('E25').Value = '  -0.77%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.518'
$ws.Range('E26').Value = '  -0.65%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.198'
$ws.Range('E27').Value = '  +3.23%  '
$ws.Range('E28').Value = '  -3.62%  '
$ws.Range('E29').Value = '  -0.69%  '
$ws.Range('E30').Value = '  +0.08%  '
$ws.Range('E31').Value = '  -0.69%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '22.46'
$ws.Range('E32').Value = '  -1.63%  '
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('E34').Value = '  -1.81%  '
$ws.Range('E35').Value = '  -3.04%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.21'
$ws.Range('E36').Value = '  -2.85%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '159.20'
$ws.Range('E37').Value = '  -0.42%  '
$ws.Range('E38').Value = '  -2.50%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.83'
$ws.Range('E39').Value = '  -1.08%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '26.86'
$ws.Range('E40').Value = '  -1.59%  '
$ws.Range('D41').Value = '2.883.77'
$ws.Range('E41').Value = '  +1.66%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.770'
$ws.Range('E42').Value = '  -2.60%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.34'
$ws.Range('E43').Value = '  -2.24%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '40.27'
$ws.Range('E44').Value = '  +0.31%  '
$ws.Range('E45').Value = '  -1.27%  '
$ws.Range('E46').Value = '  -2.78%  '
$ws.Range('E47').Value = '  -1.69%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.40'
$ws.Range('E48').Value = '  -4.43%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '317.10'
$ws.Range('E49').Value = '  -2.32%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0272'
$ws.Range('E50').Value = '  -0.35%  '
$ws.Range('E51').Value = '  +3.90%  '

$textForceCells = @('D5','D6','D7','D15','D17','D18','D19','D21','D22','D24','D26','D27','D32','D36','D37','D39','D40','D42','D43','D44','D48','D49','D50')
foreach ($c in $textForceCells) { $ws.Range($c).Style = 'Normal' }
